# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 3.56341032713086)
    3 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 0, 14.05633640148523)
    4 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    5 = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 10.19245300693656, 0, 16.0466451790395)
    6 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
